# ajout systeme de point variant selon le chap
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 5 with the "chap 2" variant scoring data (mirrors row 6) ---
$ws.Range("B5").Value = "Dechiffrer le message"
$ws.Range("D5").Value = 0

$ws.Range("F5").Value = "b"
$ws.Range("G5").Value = "g"
$ws.Range("H5").Value = "m"
$ws.Range("I5").Value = "c"
$ws.Range("J5").Value = "h"
$ws.Range("K5").Value = "n"
$ws.Range("L5").Value = "Sers toi de ce que tu as vu sur le mur pour comprendre comment faire`na -> b`nb -> c`n…`nz -> a`nIl semblerait que le message n'utilise que des lettres en minuscules"
$ws.Range("L5").WrapText = $true

$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = "b"
$ws.Range("S5").Value = "c"
$ws.Range("T5").Value = "w"
$ws.Range("U5").Value = "x"
$ws.Range("V5").Value = "z"
$ws.Range("W5").Value = "a"

$ws.Rows("5").RowHeight = 259.2

# --- Update the active sheet view / selection ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("P5").Select()
